$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").EntireColumn.Insert()

$ws.Range("F2").Value = "Address"

$ws.Range("F3").Value = "G H S SSwandenahalli"
$ws.Range("F4").Value = "G H S H L KshetraTiptur"
$ws.Range("F5").Value = "Govt. Junior CollegeBellavi"
$ws.Range("F6").Value = "Kenkere High School C N Halli"
$ws.Range("F7").Value = "G P U CollegeH S Turuvekere"
$ws.Range("F8").Value = "Govt Girls High SchoolCheluru Gubbi"
$ws.Range("F9").Value = "G H S Hanumanthapura"
$ws.Range("F10").Value = "G H P SchoolDoddagoraghattaTuruvekere"
$ws.Range("F11").Value = "U G H P SHethenahalli"
$ws.Range("F12").Value = "G H SGuleharive"
$ws.Range("F13").Value = "G H S C N Halli"
$ws.Range("F14").Value = "G J C Borana KaniveC N Hally"
$ws.Range("F15").Value = "G H S BedathurMadhugiri"
$ws.Range("F16").Value = "G H SGowdanakatteTiptur"
$ws.Range("F17").Value = "S R R H SRamanahalliC N Halli"
$ws.Range("F18").Value = "G H SKempanahalliKunigal"
$ws.Range("F19").Value = "G J C Bellavi"
$ws.Range("F20").Value = "Bhoruka High SchoolShahapur"
$ws.Range("F21").Value = "G H SHosurKunigal"
$ws.Range("F22").Value = "S V P J CollegeB H Road"
$ws.Range("F24").Value = "G H S C N Hall"
$ws.Range("F25").Value = "G H S D HosahalliKunigal (t)"
$ws.Range("F26").Value = "S V H S ThandagaTuruvekere"
$ws.Range("F27").Value = "Vidyaranya H S BommenahalliC N Halli"
$ws.Range("F28").Value = "G M H P SKyathasandva"
$ws.Range("F29").Value = "G P U C HuliyarKenkereC N Halli"
$ws.Range("F30").Value = "G J C Bellavi"
$ws.Range("F31").Value = "G H S H TammadihallyC N Hally"
$ws.Range("F32").Value = "G H SDurgadahally"
$ws.Range("F33").Value = "G H S KamalapuraC N Halli"
$ws.Range("F34").Value = "G J C AmruthurKunigal"
$ws.Range("F35").Value = "G H SP H Colony"
$ws.Range("F36").Value = "S G R H S AnekereTuruvekere"
$ws.Range("F37").Value = "G H P S KondliGubbi"
$ws.Range("F38").Value = "S C H SHarenahally GateC N Hally"
$ws.Range("F39").Value = "G H SHullenahalli"
$ws.Range("F40").Value = "Bhoruka English Medium SchoolShivapur"
$ws.Range("F41").Value = "G J C KadabaGubbi"
$ws.Range("F42").Value = "S S A H S ModuruKunigal"
$ws.Range("F43").Value = "S V A H SGollarahattiKunigal"
